$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2106918238993711
$ws.Range("C2").Value = 0.5377358490566038
$ws.Range("J2").Value = 0.01257861635220126
$ws.Range("P2").Value = 0.160377358490566
$ws.Range("S2").Value = 0.07861635220125786
$ws.Range("B3").Value = 0.02793296089385475
$ws.Range("C3").Value = 0.03910614525139665
$ws.Range("J3").Value = 0.01675977653631285
$ws.Range("P3").Value = 0.7206703910614525
$ws.Range("S3").Value = 0.1955307262569832
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.5909090909090909
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.04545454545454546
$ws.Range("D6").Value = 0.005050505050505051
$ws.Range("F6").Value = 0.02525252525252525
$ws.Range("J6").Value = 0.2575757575757576
$ws.Range("O6").Value = 0.01515151515151515
$ws.Range("Q6").Value = 0.1818181818181818
$ws.Range("R6").Value = 0.1161616161616162
$ws.Range("S6").Value = 0.3535353535353535
$ws.Range("B7").Value = 0.1151079136690648
$ws.Range("D7").Value = 0.01079136690647482
$ws.Range("E7").Value = 0.003597122302158274
$ws.Range("F7").Value = 0.04316546762589928
$ws.Range("J7").Value = 0.1618705035971223
$ws.Range("O7").Value = 0.01798561151079137
$ws.Range("Q7").Value = 0.158273381294964
$ws.Range("R7").Value = 0.06115107913669065
$ws.Range("S7").Value = 0.4280575539568345
$ws.Range("B8").Value = 0.09368191721132897
$ws.Range("D8").Value = 0.01089324618736384
$ws.Range("F8").Value = 0.04357298474945534
$ws.Range("J8").Value = 0.1132897603485839
$ws.Range("O8").Value = 0.02396514161220044
$ws.Range("Q8").Value = 0.159041394335512
$ws.Range("R8").Value = 0.1154684095860566
$ws.Range("S8").Value = 0.4400871459694989
$ws.Range("B9").Value = 0.09595959595959595
$ws.Range("D9").Value = 0.0101010101010101
$ws.Range("E9").Value = 0.005050505050505051
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.04040404040404041
$ws.Range("Q9").Value = 0.1565656565656566
$ws.Range("R9").Value = 0.1161616161616162
$ws.Range("S9").Value = 0.4090909090909091
$ws.Range("B10").Value = 0.1055389221556886
$ws.Range("D10").Value = 0.02694610778443114
$ws.Range("E10").Value = 0.0007485029940119761
$ws.Range("F10").Value = 0.06811377245508982
$ws.Range("J10").Value = 0.1407185628742515
$ws.Range("O10").Value = 0.01122754491017964
$ws.Range("Q10").Value = 0.1998502994011976
$ws.Range("R10").Value = 0.09206586826347306
$ws.Range("S10").Value = 0.3547904191616766
$ws.Range("G11").Value = 0.1642156862745098
$ws.Range("J11").Value = 0.07843137254901961
$ws.Range("K11").Value = 0.1936274509803922
$ws.Range("L11").Value = 0.5563725490196079
$ws.Range("S11").Value = 0.007352941176470588
$ws.Range("G12").Value = 0.7574468085106383
$ws.Range("J12").Value = 0.1872340425531915
$ws.Range("K12").Value = 0.02553191489361702
$ws.Range("L12").Value = 0.02127659574468085
$ws.Range("S12").Value = 0.008510638297872341
$ws.Range("F13").Value = 0.01538461538461539
$ws.Range("G13").Value = 0.6615384615384615
$ws.Range("J13").Value = 0.2769230769230769
$ws.Range("S13").Value = 0.04615384615384616
$ws.Range("F15").Value = 0.04945054945054945
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.06043956043956044
$ws.Range("J15").Value = 0.3351648351648351
$ws.Range("K15").Value = 0.09340659340659341
$ws.Range("M15").Value = 0.01098901098901099
$ws.Range("O15").Value = 0.04945054945054945
$ws.Range("S15").Value = 0.2582417582417583
$ws.Range("F16").Value = 0.009950248756218905
$ws.Range("H16").Value = 0.1592039800995025
$ws.Range("I16").Value = 0.0945273631840796
$ws.Range("J16").Value = 0.3631840796019901
$ws.Range("K16").Value = 0.1691542288557214
$ws.Range("M16").Value = 0.06965174129353234
$ws.Range("O16").Value = 0.02985074626865672
$ws.Range("S16").Value = 0.1044776119402985
$ws.Range("F17").Value = 0.0111358574610245
$ws.Range("H17").Value = 0.1915367483296214
$ws.Range("I17").Value = 0.08908685968819599
$ws.Range("J17").Value = 0.4498886414253898
$ws.Range("K17").Value = 0.1158129175946548
$ws.Range("M17").Value = 0.0289532293986637
$ws.Range("O17").Value = 0.0400890868596882
$ws.Range("S17").Value = 0.07349665924276169
$ws.Range("H18").Value = 0.1791666666666667
$ws.Range("I18").Value = 0.1083333333333333
$ws.Range("J18").Value = 0.4166666666666667
$ws.Range("K18").Value = 0.125
$ws.Range("M18").Value = 0.02083333333333333
$ws.Range("N18").Value = 0.004166666666666667
$ws.Range("O18").Value = 0.07916666666666666
$ws.Range("S18").Value = 0.06666666666666667
$ws.Range("F19").Value = 0.01305057096247961
$ws.Range("H19").Value = 0.2283849918433931
$ws.Range("I19").Value = 0.08482871125611746
$ws.Range("J19").Value = 0.3752039151712888
$ws.Range("K19").Value = 0.1411092985318108
$ws.Range("M19").Value = 0.02773246329526917
$ws.Range("N19").Value = 0.001631321370309951
$ws.Range("O19").Value = 0.05301794453507341
$ws.Range("S19").Value = 0.07504078303425775
